$d = $word.ActiveDocument
$d.Content.Find.Execute("the TCP message, and generate", $true, $false, $false, $false, $false, $true, 1, $false, "the TCP messages, filters them keeping just the interesting ones, and generate", 2)
